$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '97.630.70'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.372.70'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '252.09'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '666.19'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.43'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -6.42%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.425'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -10.81%  '
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.03'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.65%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '3.374.26'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.212'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '41.26'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '97.389.99'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.47%  '
$ws.Range('B15').Value = 'Toncoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.18'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +8.34%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000259'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -5.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.001.79'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.68'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +13.96%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.386.26'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.577'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +32.59%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.13'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.94'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +3.40%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '504.83'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.57%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.39'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -5.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0000202'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -7.71%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.30'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +3.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '95.69'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -6.38%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '12.39'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.576.00'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.70%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '11.31'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +2.67%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.996'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.191'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.60'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +22.99%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.996'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.560'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.89%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '29.00'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '7.88'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.49'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +11.90%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '534.17'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.151'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -4.56%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '24.69'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.00'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +15.16%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.853'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +4.00%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0429'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('B47').Value = 'MantraDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.69'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -4.35%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.69'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +11.07%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.69'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +10.89%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '54.22'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +6.59%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.18'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -7.21%  '
